$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Remove the _GoBack bookmark that currently sits at the end of the
#    "Last but not least is the controller..." paragraph. The engine
#    does not expose Bookmarks via the Bookmarks collection for
#    _GoBack, so we rebuild that paragraph's content (identical runs,
#    just without the trailing bookmarkStart/bookmarkEnd markers).
# -----------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Last but not least is the controller*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the controller paragraph"
}

$targetRange = $target.Range
$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00234C5B" w:rsidRDefault="00234C5B" w:rsidP="00635E0B"><w:pPr><w:pStyle w:val="Akapitzlist"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Last but not least is the controller- a singleton that connects all parts of the system together. It is illustrated by the sequence diagram</w:t></w:r><w:r w:rsidR="00F03DBB"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> of searching for lecturers. In the view the user inputs searched phrase. Then The view calls the controller, which calls the Search engine. The search engine calls different classes and at last, it returns a ready </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00F03DBB"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>ArrayList</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00F03DBB"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> of searched lecturers back to the controller, which passes it further to the view where the information is being displayed.</w:t></w:r></w:p>
'@
$targetRange.InsertXML($paraXml) | Out-Null

# -----------------------------------------------------------------
# 2) Append a brand new paragraph ("Present") at the very end of the
#    document, carrying the _GoBack bookmark that was removed above.
# -----------------------------------------------------------------
$endPos = $d.Content.End
$endRange = $d.Range($endPos, $endPos)
$newParaXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="360"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t>Present</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@
$endRange.InsertXML($newParaXml) | Out-Null
